$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlVAlignCenter - the data cells on this sheet all share one style that is
# just "vertical centered" (no custom number format / font), so re-applying
# this after touching a cell's value reproduces that same style exactly.
$xlVAlignCenter = -4108

# --- Column width changes ---
# Raw OOXML widths: F 33->32, H 33->32, J 14->34.
# COM ColumnWidth reads ~0.83 lower than the raw stored width on this sheet,
# so subtract that offset when setting to land on the exact target width.
$ws.Columns.Item(6).ColumnWidth = 31.17   # F: 33 -> 32
$ws.Columns.Item(8).ColumnWidth = 31.17   # H: 33 -> 32
$ws.Columns.Item(10).ColumnWidth = 33.17  # J: 14 -> 34

# --- Rename "CHUANG, Ivy" -> "Ivy CHUANG" in private-lesson cells ---
$ws.Range("B7").Value = "Private lesson with Ivy CHUANG"
$ws.Range("H7").Value = "Private lesson with Ivy CHUANG"
$ws.Range("F11").Value = "Private lesson with Ivy CHUANG"

# --- Replace "Free Time" with the new Master class activity ---
$ws.Range("J7").Value = "Master class with Ivy & Stephane"
$ws.Range("J20").Value = "Master class with Ivy & Stephane"

# --- Split the J20:J27 merged block: J20:J23 keeps the master class, and
#     J24:J27 becomes its own separate "Free Time" block (new J24 cell) ---
$ws.Range("J20:J27").UnMerge()
$ws.Range("J20:J27").Clear()
$ws.Range("J20:J23").Merge()
$ws.Range("J20").Value = "Master class with Ivy & Stephane"
$ws.Range("J20").VerticalAlignment = $xlVAlignCenter
$ws.Range("J24:J27").Merge()
$ws.Range("J24").Value = "Free Time"
$ws.Range("J24").VerticalAlignment = $xlVAlignCenter

# --- Remove the Day1-5 "Free Time" row-32 block entirely (activities after
#     17:00 are no longer shown on the student timetable) ---
$ws.Range("B32:B39").UnMerge()
$ws.Range("D32:D39").UnMerge()
$ws.Range("F32:F39").UnMerge()
$ws.Range("H32:H39").UnMerge()
$ws.Range("J32:J39").UnMerge()
$ws.Range("B32").Clear()
$ws.Range("D32").Clear()
$ws.Range("F32").Clear()
$ws.Range("H32").Clear()
$ws.Range("J32").Clear()

# --- Shrink the B/D/F/H 28:31 merges down to 28:30 to match (row 31, the
#     17:00 slot, is no longer part of these activity blocks) ---
$ws.Range("B28:B31").UnMerge()
$ws.Range("B28:B31").Clear()
$ws.Range("B28:B30").Merge()
$ws.Range("B28").Value = "Acting class"
$ws.Range("B28").VerticalAlignment = $xlVAlignCenter

$ws.Range("D28:D31").UnMerge()
$ws.Range("D28:D31").Clear()
$ws.Range("D28:D30").Merge()
$ws.Range("D28").Value = "Acting class"
$ws.Range("D28").VerticalAlignment = $xlVAlignCenter

$ws.Range("F28:F31").UnMerge()
$ws.Range("F28:F31").Clear()
$ws.Range("F28:F30").Merge()
$ws.Range("F28").Value = "Acting class"
$ws.Range("F28").VerticalAlignment = $xlVAlignCenter

$ws.Range("H28:H31").UnMerge()
$ws.Range("H28:H31").Clear()
$ws.Range("H28:H30").Merge()
$ws.Range("H28").Value = "Acting class"
$ws.Range("H28").VerticalAlignment = $xlVAlignCenter

# --- Shrink J28:J31 down to J28:J30 as well ---
$ws.Range("J28:J31").UnMerge()
$ws.Range("J28:J31").Clear()
$ws.Range("J28:J30").Merge()
$ws.Range("J28").Value = "Acting class"
$ws.Range("J28").VerticalAlignment = $xlVAlignCenter
